# Update the marksheet "Corr/total marks" figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking (points per correct answer): 3 -> 5
$ws.Range("B11").Value = 5

# Total correct marks: 60 -> 100
$ws.Range("B12").Value = 100

# Total / max marks label: 60/84 -> 100/140
$ws.Range("E12").Value = "100/140"
